$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "release/8.0.14"
$ws.Range("B17").Value = "X"
$ws.Range("C17").Value = "X"
$ws.Range("D17").Value = "X"
$ws.Range("E17").Value = "X"
